$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 289 (shifts existing rows 289-311 down to 290-312,
# and extends the sheet dimension to R312 automatically).
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with the new weekly price record.
$ws.Range("A289").Value = 3
$ws.Range("B289").Value = "Femacal de La Calera"
$ws.Range("C289").Value = "Coquimbo"
$ws.Range("D289").Value = "2022-07-04"
$ws.Range("E289").Value = 5
$ws.Range("F289").Value = 100112001
$ws.Range("G289").Value = "Berenjena"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 60
$ws.Range("K289").Value = 7000
$ws.Range("L289").Value = 7000
$ws.Range("M289").Value = 7000
$ws.Range("N289").Value = "`$/caja 60 unidades"
$ws.Range("O289").Value = "Región de Arica y Parinacota"
$ws.Range("P289").Value = 117
$ws.Range("Q289").Value = 60
$ws.Range("R289").Value = "Hortaliza"
